$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2, 4) '27.504.12'
Set-TextValue $ws.Cells.Item(2, 5) '  +6.82%  '

Set-TextValue $ws.Cells.Item(3, 4) '1.814.23'
Set-TextValue $ws.Cells.Item(3, 5) '  +6.55%  '

Set-TextValue $ws.Cells.Item(4, 5) '  -0.24%  '

Set-TextValue $ws.Cells.Item(5, 4) '345.37'
Set-TextValue $ws.Cells.Item(5, 5) '  +4.50%  '

Set-TextValue $ws.Cells.Item(6, 4) '0.9986'
Set-TextValue $ws.Cells.Item(6, 5) '  -0.07%  '

Set-TextValue $ws.Cells.Item(7, 5) '  +4.28%  '

Set-TextValue $ws.Cells.Item(8, 2) 'Cardano'
Set-TextValue $ws.Cells.Item(8, 3) 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Cells.Item(8, 4) '0.3526'
Set-TextValue $ws.Cells.Item(8, 5) '  +6.78%  '

Set-TextValue $ws.Cells.Item(9, 2) 'OKB'
Set-TextValue $ws.Cells.Item(9, 3) 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws.Cells.Item(9, 4) '50.08'
Set-TextValue $ws.Cells.Item(9, 5) '  +3.23%  '

Set-TextValue $ws.Cells.Item(10, 4) '1.238'
Set-TextValue $ws.Cells.Item(10, 5) '  +6.11%  '

Set-TextValue $ws.Cells.Item(11, 4) '0.07754'
Set-TextValue $ws.Cells.Item(11, 5) '  +5.90%  '

Set-TextValue $ws.Cells.Item(12, 4) '1.000'
Set-TextValue $ws.Cells.Item(12, 5) '  -0.12%  '

Set-TextValue $ws.Cells.Item(13, 4) '22.56'
Set-TextValue $ws.Cells.Item(13, 5) '  +12.92%  '

Set-TextValue $ws.Cells.Item(14, 4) '6.639'
Set-TextValue $ws.Cells.Item(14, 5) '  +7.37%  '

Set-TextValue $ws.Cells.Item(15, 4) '7.230'
Set-TextValue $ws.Cells.Item(15, 5) '  +5.62%  '

Set-TextValue $ws.Cells.Item(16, 4) '1.813.38'
Set-TextValue $ws.Cells.Item(16, 5) '  +6.56%  '

Set-TextValue $ws.Cells.Item(17, 5) '  +6.10%  '

Set-TextValue $ws.Cells.Item(18, 4) '0.06739'
Set-TextValue $ws.Cells.Item(18, 5) '  +1.68%  '

Set-TextValue $ws.Cells.Item(19, 4) '86.99'
Set-TextValue $ws.Cells.Item(19, 5) '  +7.35%  '

Set-TextValue $ws.Cells.Item(20, 4) '0.9988'
Set-TextValue $ws.Cells.Item(20, 5) '  -0.11%  '

Set-TextValue $ws.Cells.Item(21, 4) '17.83'
Set-TextValue $ws.Cells.Item(21, 5) '  +10.62%  '

Set-TextValue $ws.Cells.Item(22, 4) '6.543'
Set-TextValue $ws.Cells.Item(22, 5) '  +8.51%  '

Set-TextValue $ws.Cells.Item(23, 4) '13.19'
Set-TextValue $ws.Cells.Item(23, 5) '  +1.88%  '

Set-TextValue $ws.Cells.Item(24, 4) '27.493.67'

Set-TextValue $ws.Cells.Item(25, 4) '2.469'
Set-TextValue $ws.Cells.Item(25, 5) '  +0.47%  '

Set-TextValue $ws.Cells.Item(26, 4) '2.693'
Set-TextValue $ws.Cells.Item(26, 5) '  +8.94%  '

Set-TextValue $ws.Cells.Item(27, 4) '22.12'
Set-TextValue $ws.Cells.Item(27, 5) '  +15.42%  '

Set-TextValue $ws.Cells.Item(28, 4) '1.506'
Set-TextValue $ws.Cells.Item(28, 5) '  +17.07%  '

Set-TextValue $ws.Cells.Item(29, 4) '154.25'
Set-TextValue $ws.Cells.Item(29, 5) '  +3.20%  '

Set-TextValue $ws.Cells.Item(30, 4) '2.015.23'
Set-TextValue $ws.Cells.Item(30, 5) '  +6.52%  '

Set-TextValue $ws.Cells.Item(31, 4) '137.01'
Set-TextValue $ws.Cells.Item(31, 5) '  +7.01%  '

Set-TextValue $ws.Cells.Item(32, 4) '6.410'
Set-TextValue $ws.Cells.Item(32, 5) '  +8.43%  '

Set-TextValue $ws.Cells.Item(33, 4) '4.086'
Set-TextValue $ws.Cells.Item(33, 5) '  -0.27%  '

Set-TextValue $ws.Cells.Item(34, 4) '13.97'
Set-TextValue $ws.Cells.Item(34, 5) '  +9.02%  '

Set-TextValue $ws.Cells.Item(35, 4) '0.08838'
Set-TextValue $ws.Cells.Item(35, 5) '  +4.29%  '

Set-TextValue $ws.Cells.Item(36, 4) '1.726'
Set-TextValue $ws.Cells.Item(36, 5) '  +0.63%  '

Set-TextValue $ws.Cells.Item(37, 4) '5.683'
Set-TextValue $ws.Cells.Item(37, 5) '  +6.96%  '

Set-TextValue $ws.Cells.Item(38, 4) '0.7104'
Set-TextValue $ws.Cells.Item(38, 5) '  +16.58%  '

Set-TextValue $ws.Cells.Item(39, 4) '0.06580'
Set-TextValue $ws.Cells.Item(39, 5) '  +6.62%  '

Set-TextValue $ws.Cells.Item(40, 4) '0.02430'
Set-TextValue $ws.Cells.Item(40, 5) '  +8.52%  '

Set-TextValue $ws.Cells.Item(41, 4) '0.2274'
Set-TextValue $ws.Cells.Item(41, 5) '  +7.72%  '

Set-TextValue $ws.Cells.Item(42, 4) '9.025'
Set-TextValue $ws.Cells.Item(42, 5) '  +6.27%  '

Set-TextValue $ws.Cells.Item(43, 4) '1.298'
Set-TextValue $ws.Cells.Item(43, 5) '  +1.98%  '

Set-TextValue $ws.Cells.Item(44, 4) '15.06'
Set-TextValue $ws.Cells.Item(44, 5) '  +2.49%  '

Set-TextValue $ws.Cells.Item(45, 4) '0.6623'
Set-TextValue $ws.Cells.Item(45, 5) '  +13.91%  '

Set-TextValue $ws.Cells.Item(46, 4) '0.9980'

Set-TextValue $ws.Cells.Item(47, 4) '4.048'
Set-TextValue $ws.Cells.Item(47, 5) '  +5.55%  '

Set-TextValue $ws.Cells.Item(48, 4) '2.192'
Set-TextValue $ws.Cells.Item(48, 5) '  +9.95%  '

Set-TextValue $ws.Cells.Item(49, 4) '133.31'
Set-TextValue $ws.Cells.Item(49, 5) '  +5.52%  '

Set-TextValue $ws.Cells.Item(50, 4) '0.07366'
Set-TextValue $ws.Cells.Item(50, 5) '  +2.12%  '

Set-TextValue $ws.Cells.Item(51, 4) '80.89'
Set-TextValue $ws.Cells.Item(51, 5) '  +6.01%  '

Write-Host "Applied 99 cell updates"